$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.298.06'
$ws.Range("E2").Value = '  +1.95%  '
$ws.Range("D3").Value = '1.915.30'
$ws.Range("E3").Value = '  +2.16%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.ClearFormats()
$ws.Range("E4").Value = '  -0.36%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '327.96'
$c.ClearFormats()
$ws.Range("E5").Value = '  +0.85%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.ClearFormats()
$ws.Range("E6").Value = '  -0.21%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4619'
$c.ClearFormats()
$ws.Range("E7").Value = '  +0.58%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3952'
$c.ClearFormats()
$ws.Range("E8").Value = '  +2.37%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '46.78'
$c.ClearFormats()
$ws.Range("E9").Value = '  +1.68%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.07954'
$c.ClearFormats()
$ws.Range("E10").Value = '  +1.25%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.ClearFormats()
$ws.Range("E11").Value = '  +0.90%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '22.41'
$c.ClearFormats()
$ws.Range("E12").Value = '  +2.87%  '
$ws.Range("D13").Value = '1.913.76'
$ws.Range("E13").Value = '  +1.58%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '7.107'
$c.ClearFormats()
$ws.Range("E14").Value = '  +1.66%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.781'
$c.ClearFormats()
$ws.Range("E15").Value = '  +1.20%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.06954'
$c.ClearFormats()
$ws.Range("E16").Value = '  +0.06%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '88.58'
$c.ClearFormats()
$ws.Range("E17").Value = '  +0.16%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.ClearFormats()
$ws.Range("E18").Value = '  -0.33%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.00001009'
$c.ClearFormats()
$ws.Range("E19").Value = '  +0.41%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '17.17'
$c.ClearFormats()
$ws.Range("E20").Value = '  +1.80%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("D22").Value = '29.272.61'
$ws.Range("E22").Value = '  +1.79%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.371'
$c.ClearFormats()
$ws.Range("E23").Value = '  +1.78%  '
$ws.Range("E24").Value = '  +0.92%  '
$ws.Range("D25").Value = '2.119.23'
$ws.Range("E25").Value = '  +0.50%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.064'
$c.ClearFormats()
$ws.Range("E26").Value = '  -3.25%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '156.70'
$c.ClearFormats()
$ws.Range("E27").Value = '  +2.15%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '19.52'
$c.ClearFormats()
$ws.Range("E28").Value = '  +1.49%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '6.098'
$c.ClearFormats()
$ws.Range("E29").Value = '  +5.69%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.002'
$c.ClearFormats()
$ws.Range("E30").Value = '  +1.79%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '119.06'
$c.ClearFormats()
$ws.Range("E31").Value = '  +0.06%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.09397'
$c.ClearFormats()
$ws.Range("E32").Value = '  +0.77%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.9285'
$c.ClearFormats()
$ws.Range("E33").Value = '  +1.04%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.353'
$c.ClearFormats()
$ws.Range("E34").Value = '  +0.98%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.362'
$c.ClearFormats()
$ws.Range("E35").Value = '  +1.65%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '3.275'
$c.ClearFormats()
$ws.Range("E36").Value = '  -1.47%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.210'
$c.ClearFormats()
$ws.Range("E37").Value = '  +4.78%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.05848'
$c.ClearFormats()
$ws.Range("E38").Value = '  +1.44%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.02110'
$c.ClearFormats()
$ws.Range("E39").Value = '  +1.89%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '7.968'
$c.ClearFormats()
$ws.Range("E40").Value = '  +3.39%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range("E41").Value = '  -0.23%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.5761'
$c.ClearFormats()
$ws.Range("E42").Value = '  +2.15%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.1805'
$c.ClearFormats()
$ws.Range("E43").Value = '  +0.86%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '9.999'
$c.ClearFormats()
$ws.Range("E44").Value = '  +0.93%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.272'
$c.ClearFormats()
$ws.Range("E45").Value = '  +6.06%  '
$ws.Range("E46").Value = '  +1.63%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.5429'
$c.ClearFormats()
$ws.Range("E47").Value = '  +2.50%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.07080'
$c.ClearFormats()
$ws.Range("E48").Value = '  -1.85%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.880'
$c.ClearFormats()
$ws.Range("E49").Value = '  +3.02%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.565'
$c.ClearFormats()
$ws.Range("E50").Value = '  +6.33%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '113.20'
$c.ClearFormats()
$ws.Range("E51").Value = '  -0.35%  '
